{"js": "// Replace the worksheet date and every \"a\u00d7b=c\" multiplication answer with\n// its updated value, per the commit diff (2024-01-07 -> 2024-01-08, and\n// 25 multiplication results recomputed for the new date).\nconst replacements = [\n  [\"2024-01-07 Sunday\", \"2024-01-08 Monday\"],\n  [\"58\u00d740=2320\", \"69\u00d783=5727\"],\n  [\"71\u00d720=1420\", \"23\u00d799=2277\"],\n  [\"78\u00d742=3276\", \"53\u00d785=4505\"],\n  [\"40\u00d747=1880\", \"53\u00d794=4982\"],\n  [\"54\u00d718=972\", \"90\u00d758=5220\"],\n  [\"33\u00d767=2211\", \"68\u00d767=4556\"],\n  [\"70\u00d712=840\", \"87\u00d717=1479\"],\n  [\"28\u00d717=476\", \"84\u00d799=8316\"],\n  [\"79\u00d766=5214\", \"95\u00d714=1330\"],\n  [\"64\u00d759=3776\", \"38\u00d773=2774\"],\n  [\"97\u00d756=5432\", \"32\u00d774=2368\"],\n  [\"88\u00d725=2200\", \"61\u00d717=1037\"],\n  [\"30\u00d715=450\", \"87\u00d782=7134\"],\n  [\"70\u00d730=2100\", \"32\u00d750=1600\"],\n  [\"20\u00d730=600\", \"99\u00d779=7821\"],\n  [\"71\u00d749=3479\", \"56\u00d748=2688\"],\n  [\"16\u00d781=1296\", \"59\u00d735=2065\"],\n  [\"54\u00d738=2052\", \"60\u00d717=1020\"],\n  [\"84\u00d782=6888\", \"37\u00d745=1665\"],\n  [\"96\u00d723=2208\", \"59\u00d717=1003\"],\n  [\"16\u00d750=800\", \"22\u00d725=550\"],\n  [\"58\u00d714=812\", \"31\u00d761=1891\"],\n  [\"62\u00d747=2914\", \"50\u00d728=1400\"],\n  [\"26\u00d790=2340\", \"98\u00d769=6762\"],\n  [\"34\u00d726=884\", \"85\u00d750=4250\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the worksheet date and every \"a\u00d7b=c\" multiplication answer with\n# its updated value, per the commit diff (2024-01-07 -> 2024-01-08, and\n# 25 multiplication results recomputed for the new date).\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2024-01-07 Sunday\", \"2024-01-08 Monday\"),\n    @(\"58\u00d740=2320\", \"69\u00d783=5727\"),\n    @(\"71\u00d720=1420\", \"23\u00d799=2277\"),\n    @(\"78\u00d742=3276\", \"53\u00d785=4505\"),\n    @(\"40\u00d747=1880\", \"53\u00d794=4982\"),\n    @(\"54\u00d718=972\", \"90\u00d758=5220\"),\n    @(\"33\u00d767=2211\", \"68\u00d767=4556\"),\n    @(\"70\u00d712=840\", \"87\u00d717=1479\"),\n    @(\"28\u00d717=476\", \"84\u00d799=8316\"),\n    @(\"79\u00d766=5214\", \"95\u00d714=1330\"),\n    @(\"64\u00d759=3776\", \"38\u00d773=2774\"),\n    @(\"97\u00d756=5432\", \"32\u00d774=2368\"),\n    @(\"88\u00d725=2200\", \"61\u00d717=1037\"),\n    @(\"30\u00d715=450\", \"87\u00d782=7134\"),\n    @(\"70\u00d730=2100\", \"32\u00d750=1600\"),\n    @(\"20\u00d730=600\", \"99\u00d779=7821\"),\n    @(\"71\u00d749=3479\", \"56\u00d748=2688\"),\n    @(\"16\u00d781=1296\", \"59\u00d735=2065\"),\n    @(\"54\u00d738=2052\", \"60\u00d717=1020\"),\n    @(\"84\u00d782=6888\", \"37\u00d745=1665\"),\n    @(\"96\u00d723=2208\", \"59\u00d717=1003\"),\n    @(\"16\u00d750=800\", \"22\u00d725=550\"),\n    @(\"58\u00d714=812\", \"31\u00d761=1891\"),\n    @(\"62\u00d747=2914\", \"50\u00d728=1400\"),\n    @(\"26\u00d790=2340\", \"98\u00d769=6762\"),\n    @(\"34\u00d726=884\", \"85\u00d750=4250\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 1  # wdFindContinue\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null\n}\n"}
